$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.9702141024943302
$ws1.Range("C2").Value = -0.4881301422590007
$ws1.Range("B3").Value = -0.1861873666424456
$ws1.Range("C3").Value = -0.9090792945492638
$ws1.Range("B4").Value = 0.8606357959804208
$ws1.Range("C4").Value = -0.6678617948341953

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -0.9843908645868247
$ws2.Range("C2").Value = -0.4004554809490697
$ws2.Range("B3").Value = 0.5037057425529781
$ws2.Range("C3").Value = 0.7293369576195338
$ws2.Range("B4").Value = -1.118079193513543
$ws2.Range("C4").Value = 0.6311520462704328
